# plusieurs candidatures au 27/08/2018
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Helper: write a date into column A re-using the existing date style (s="1")
# by copying the format from an already-dated cell (A16) instead of letting
# the engine mint a brand-new number-format style for every new cell.
function Set-DateCell($row, $serial) {
    $ws.Cells.Item(16, 1).Copy() | Out-Null
    $ws.Cells.Item($row, 1).PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($row, 1).Value = $serial
}

$dateSerial = 43339

# Row 17 - pôle emploi, Développeur Intégrateur de Logiciels, Villeneuve-d'Ascq
Set-DateCell 17 $dateSerial
$ws.Cells.Item(17, 3).Value = "Développeur Intégrateur de Logiciels"
$ws.Cells.Item(17, 4).Value = "CDI"
$ws.Cells.Item(17, 5).Value = "Temps Plein"
$ws.Cells.Item(17, 6).Value = "Villeneuve-d'Ascq"
$ws.Cells.Item(17, 7).Value = "pôle emploi"

# Row 18 - Supplay via NordJob, Développeur Web, Arras
Set-DateCell 18 $dateSerial
$ws.Cells.Item(18, 3).Value = "Développeur Web"
$ws.Cells.Item(18, 6).Value = "Arras"
$ws.Cells.Item(18, 7).Value = "Supplay via NordJob"

# Row 19 - Développeur Web, Lille, Monster (entreprise filled in afterwards)
Set-DateCell 19 $dateSerial
$ws.Cells.Item(19, 3).Value = "Développeur Web"
$ws.Cells.Item(19, 4).Value = "CDI"
$ws.Cells.Item(19, 5).Value = "Temps Plein"
$ws.Cells.Item(19, 6).Value = "Lille"
$ws.Cells.Item(19, 7).Value = "Monster"

# Row 20 - Kalyptus, Développeur PHP, CDi, Villeneuve-d'Ascq, Monster
Set-DateCell 20 $dateSerial
$ws.Cells.Item(20, 2).Value = "Kalyptus"
$ws.Cells.Item(20, 3).Value = "Développeur PHP"
$ws.Cells.Item(20, 4).Value = "CDi"
$ws.Cells.Item(20, 5).Value = "Temps Plein"
$ws.Cells.Item(20, 6).Value = "Villeneuve-d'Ascq"
$ws.Cells.Item(20, 7).Value = "Monster"

# Entreprise for row 19 entered as an afterthought, after row 20 was filled in
$ws.Cells.Item(19, 2).Value = "Studio RH"

$ws.Range("A21").Select() | Out-Null
